# Sum expenses + save cost as float
#
# The "2019" expense sheet used to hold three string-typed cost rows
# (Beef Steak/350, Candy/20, Kra Pow Kai/55). Replace that with two rows
# whose costs are real numbers, add a Cha Yen row, drop the Candy /
# Kra Pow Kai rows, and total the costs with a SUM formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# Row 1: Cha Yen, cost 45 (numeric, not text)
$ws.Range("B1").Value = "Cha Yen"
$ws.Range("C1").Value = 45

# Row 2: Beef Steak, cost 350 (numeric, not text)
$ws.Range("B2").Value = "Beef Steak"
$ws.Range("C2").Value = 350

# Row 3 (old "Kra Pow Kai" row) is no longer needed
$ws.Rows.Item(3).Clear()

# Total of the expenses in column C
$ws.Range("D1").Formula = "=SUM(C1:C2)"
